$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")
$ws.Range("B22").Value = 3
